$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:126 down to 44:127.
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new weekly price record (same dimension
# columns as the surrounding "Camote" rows, new date + prices + origin).
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(43, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(43, 4).Value = 45281
$ws.Cells.Item(43, 5).Value = 15
$ws.Cells.Item(43, 6).Value = 100112045
$ws.Cells.Item(43, 7).Value = "Zapallo"
$ws.Cells.Item(43, 8).Value = "Camote"
$ws.Cells.Item(43, 9).Value = "1a nueva(o)"
$ws.Cells.Item(43, 10).Value = 1000
$ws.Cells.Item(43, 11).Value = 1400
$ws.Cells.Item(43, 12).Value = 1500
$ws.Cells.Item(43, 13).Value = 1450
$ws.Cells.Item(43, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(43, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(43, 16).Value = 1450
$ws.Cells.Item(43, 17).Value = 1
$ws.Cells.Item(43, 18).Value = "Hortaliza"
